$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (6) - this shifts Data and everything
# after it one column to the right, making room for the new "Ano" field.
$ws.Columns.Item(6).Insert()

# Set the new column's header to "Ano"
$ws.Range("F1").Value = "Ano"

# New column inherits the custom width of the neighboring "Procedencia" column
# (matches native Excel insert-column width-carryover behavior).
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# Move the active selection to F1 (matches the post-edit selection in the file)
$ws.Range("F1").Select()
